$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text change: "Sale Price" -> "price" ---
$ws.Range("B1").Value = "price"

# --- Number formats ---
# New currency format for the price column (B), applied to the data rows.
$ws.Range("B2:B3").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"

# "Small Tag Quantity" column (C) becomes an integer quantity format.
$ws.Range("C1:C3").NumberFormat = "0"

# --- Data edits ---
$ws.Range("A2").Value = "273-1126"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

$ws.Range("A3").Value = 2700217
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 4
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# --- Selection / view state ---
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollColumn = 1
